$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8365
    3  = 7878
    8  = 131
    9  = 125
    10 = 174
    12 = 712
    13 = 131
    14 = 1867
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
